$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column M (Dezembro) values that were previously 0
$ws.Range("M2").Value = 22862.26
$ws.Range("M3").Value = 8124.5
$ws.Range("M4").Value = 5434.87
$ws.Range("M5").Value = 2422
$ws.Range("M6").Value = 38843.63

# Update column AG (total) values to reflect the updated row totals
$ws.Range("AG2").Value = 148378.57
$ws.Range("AG3").Value = 58579.9
$ws.Range("AG4").Value = 40219.77
$ws.Range("AG5").Value = 35361.2
$ws.Range("AG6").Value = 282539.44
